$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers (column set expanded from A:AE to A:AX; several
#     metrics renamed/reordered and new weight/alfa/sigma/gini/median
#     columns + three summary columns (Sample_entropy, info_entropy,
#     PLZC) added) ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Samples"
$ws.Range("C1").Value = "Samples_worn"
$ws.Range("D1").Value = "Epochs of 1 minute"
$ws.Range("E1").Value = "average_activity_level"
$ws.Range("F1").Value = "std_change"
$ws.Range("G1").Value = "per_change"
$ws.Range("H1").Value = "sedentairy_count"
$ws.Range("I1").Value = "sedentairy_perc"
$ws.Range("J1").Value = "epochs_sedentairy_perc"
$ws.Range("K1").Value = "weight_median_sedentairy"
$ws.Range("L1").Value = "alfa_sedentairy"
$ws.Range("M1").Value = "sigma_sedentairy"
$ws.Range("N1").Value = "gini_sedentairy"
$ws.Range("O1").Value = "sedentairy_median_length"
$ws.Range("P1").Value = "epochs_sedentairy_average_length"
$ws.Range("Q1").Value = "epochs_sedentairy_max_length"
$ws.Range("R1").Value = "light_count"
$ws.Range("S1").Value = "light_perc"
$ws.Range("T1").Value = "epochs_light_perc"
$ws.Range("U1").Value = "weight_median_light"
$ws.Range("V1").Value = "alfa_light"
$ws.Range("W1").Value = "sigma_light"
$ws.Range("X1").Value = "gini_light"
$ws.Range("Y1").Value = "light_median_length"
$ws.Range("Z1").Value = "epochs_light_average_length"
$ws.Range("AA1").Value = "epochs_light_max_length"
$ws.Range("AB1").Value = "moderate_count"
$ws.Range("AC1").Value = "moderate_perc"
$ws.Range("AD1").Value = "epochs_moderate_perc"
$ws.Range("AE1").Value = "weight_median_moderate"
$ws.Range("AF1").Value = "alfa_moderate"
$ws.Range("AG1").Value = "sigma_moderate"
$ws.Range("AH1").Value = "gini_moderate"
$ws.Range("AI1").Value = "moderate_median_length"
$ws.Range("AJ1").Value = "epochs_moderate_average_length"
$ws.Range("AK1").Value = "epochs_moderate_max_length"
$ws.Range("AL1").Value = "vigorous_count"
$ws.Range("AM1").Value = "vigorous_perc"
$ws.Range("AN1").Value = "epochs_vigorous_perc"
$ws.Range("AO1").Value = "weight_median_vigorous"
$ws.Range("AP1").Value = "alfa_vigorous"
$ws.Range("AQ1").Value = "sigma_vigorous"
$ws.Range("AR1").Value = "gini_vigorous"
$ws.Range("AS1").Value = "vigorous_median_length"
$ws.Range("AT1").Value = "epochs_vigorous_average_length"
$ws.Range("AU1").Value = "epochs_vigorous_max_length"
$ws.Range("AV1").Value = "Sample_entropy"
$ws.Range("AW1").Value = "info_entropy"
$ws.Range("AX1").Value = "PLZC"

# New header cells (K1:AX1 is beyond the former A1:AE1 header range,
# or were newly created after the shift) need the bold/centered/
# bordered header style -- copy formatting from an existing header
# cell (style index 1) instead of re-deriving it by hand.
$ws.Range("J1").Copy()
$ws.Range("K1:AX1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2 values ---
$ws.Range("A2").Value = "OBS.002.csv"
$ws.Range("B2").Value = 179730
$ws.Range("C2").Value = 179730
$ws.Range("D2").Value = 239
$ws.Range("E2").Value = 1.038
$ws.Range("F2").Value = 0.583
$ws.Range("G2").Value = 33.891
$ws.Range("H2").Value = 69
$ws.Range("I2").Value = 28.87
$ws.Range("J2").Value = 37.805
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 2.606
$ws.Range("M2").Value = 0.288
$ws.Range("N2").Value = 0.359
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 2.226
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 125
$ws.Range("S2").Value = 52.301
$ws.Range("T2").Value = 41.463
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 2.175
$ws.Range("W2").Value = 0.202
$ws.Range("X2").Value = 0.509
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 3.676
$ws.Range("AA2").Value = 16
$ws.Range("AB2").Value = 12
$ws.Range("AC2").Value = 5.021
$ws.Range("AD2").Value = 12.195
$ws.Range("AE2").Value = 1
$ws.Range("AF2").Value = 8.213
$ws.Range("AG2").Value = 2.281
$ws.Range("AH2").Value = 0.133
$ws.Range("AI2").Value = 1
$ws.Range("AJ2").Value = 1.2
$ws.Range("AK2").Value = 2
$ws.Range("AL2").Value = 33
$ws.Range("AM2").Value = 13.808
$ws.Range("AN2").Value = 8.537
$ws.Range("AO2").Value = 5
$ws.Range("AP2").Value = 1.672
$ws.Range("AQ2").Value = 0.254
$ws.Range("AR2").Value = 0.313
$ws.Range("AS2").Value = 5
$ws.Range("AT2").Value = 4.714
$ws.Range("AU2").Value = 7
$ws.Range("AV2").Value = 0.507
$ws.Range("AW2").Value = 5.375
$ws.Range("AX2").Value = 0.697
